## Insert a new row (row 34) defining "measurement datum" (IAO:0000109),
## which shifts the existing rows 34-56 down to 35-57 and extends the
## used range from A1:V56 to A1:V57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 34 ("minimum dose of
# pharmacological substance"), pushing it (and everything below) down.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new term.
$ws.Cells.Item(34, 1).Value  = "IAO:0000109"
$ws.Cells.Item(34, 2).Value  = "measurement datum"
$ws.Cells.Item(34, 3).Value  = "A measurement datum is an information content entity that is a recording of the output of a measurement such as produced by a device."
$ws.Cells.Item(34, 4).Value  = "information content entity"
$ws.Cells.Item(34, 16).Value = "LSR 1"
$ws.Cells.Item(34, 17).Value = "Intervention content and delivery"
$ws.Cells.Item(34, 19).Value = "External"
$ws.Cells.Item(34, 22).Value = "PS"
